# Assets/DataTable/Data.xlsx : auto-generate an "Enums" lookup sheet,
# move it to the front of the workbook and expose it via a defined name
# ("EnumGroup") so other sheets / game code can reference the group
# column directly.   ("feat : emum 자동생성")

$wb = $excel.ActiveWorkbook

# --- 1. Move "Enums" to be the first sheet in the workbook --------------
$enums = $wb.Worksheets.Item("Enums")
$enums.Move($wb.Worksheets.Item(1)) | Out-Null

# Re-fetch a fresh reference to the (now relocated) sheet - after Move()
# the old COM reference can end up pointing at stale/duplicated state.
$enums = $wb.Worksheets.Item("Enums")

# --- 2. Populate the "Enums" sheet header / type rows --------------------
$enums.Range("A1").Value = "EnumGroup"
$enums.Range("B1").Value = "ID"
$enums.Range("C1").Value = "Name"
$enums.Range("D1").Value = "Description"

$enums.Range("A2").Value = "string"
$enums.Range("B2").Value = "int"
$enums.Range("C2").Value = "string"
$enums.Range("D2").Value = "string"

# Column widths (~16 and ~12.75 characters).
$enums.Columns.Item(1).ColumnWidth = 15.285714285714286
$enums.Columns.Item(3).ColumnWidth = 12

# --- 3. Make "Enums" the active sheet / tab ------------------------------
$enums.Activate() | Out-Null
$enums.Range("A3:D4").Select() | Out-Null

# --- 4. Register the EnumGroup defined name ------------------------------
$wb.Names.Add("EnumGroup", "=Enums!`$A:`$A") | Out-Null

# --- 5. Restore the other sheets' (cosmetic) active-cell selections ------
$ws = $wb.Worksheets.Item("EnemyData")
$ws.Activate() | Out-Null
$ws.Range("U10").Select() | Out-Null

$ws = $wb.Worksheets.Item("DestinyData")
$ws.Activate() | Out-Null
$ws.Range("P22:Q22").Select() | Out-Null

$ws = $wb.Worksheets.Item("DestinyEffectData")
$ws.Activate() | Out-Null
$ws.Range("C5").Select() | Out-Null

$ws = $wb.Worksheets.Item("ItemData")
$ws.Activate() | Out-Null
$ws.Range("R11").Select() | Out-Null

# Leave the workbook focused back on "Enums", matching the saved view.
$enums = $wb.Worksheets.Item("Enums")
$enums.Activate() | Out-Null
$enums.Range("A3:D4").Select() | Out-Null
